# Apply edit: rename sheet tab + update region values (remove years without
# data on cattle herd size) per commit message
# "ajuste no cálculo para remover anos sem dado sobre o efetivo do rebanho"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab: 2010-2025 -> 2010-2023
$ws.Name = "g3.5c Aumento(2010-2023)"

# Update region names and values (B column) for rows 2-9
$ws.Range("A2").Value = "Sergipe"
$ws.Range("B2").Value = 135.911515695629

$ws.Range("A3").Value = "Rondônia"
$ws.Range("B3").Value = 51.71609713993448

$ws.Range("A4").Value = "Mato Grosso"
$ws.Range("B4").Value = 45.7040124133387

$ws.Range("A5").Value = "Tocantins"
$ws.Range("B5").Value = 38.24931410435322

$ws.Range("A6").Value = "Pará"
$ws.Range("B6").Value = 36.70245128515432

$ws.Range("A7").Value = "Goiás"
$ws.Range("B7").Value = 35.50290489692467

$ws.Range("A8").Value = "Brasil"
$ws.Range("B8").Value = 16.47549473420316

$ws.Range("A9").Value = "Nordeste"
$ws.Range("B9").Value = -8.811494134209596
